$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps its original text formatting for purely
# numeric-looking values (Excel would otherwise coerce them to floats and
# lose trailing zeros / fixed formatting, e.g. "0.0000280" -> 2.8E-05).

$ws.Range("D2").Value = "84.455.97"
$ws.Range("E2").Value = "  +5.89%  "

$ws.Range("D3").Value = "3.285.57"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.60"
$ws.Range("E5").Value = "  +3.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "633.26"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.327"
$ws.Range("E7").Value = "  +24.17%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -2.03%  "

$ws.Range("D10").Value = "3.282.78"
$ws.Range("E10").Value = "  +2.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.603"
$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000280"
$ws.Range("E12").Value = "  +4.43%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "3.882.58"
$ws.Range("E14").Value = "  +2.13%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.91"
$ws.Range("E15").Value = "  +3.53%  "

$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  -0.40%  "

$ws.Range("D17").Value = "84.327.77"
$ws.Range("E17").Value = "  +5.99%  "

$ws.Range("D18").Value = "3.276.62"
$ws.Range("E18").Value = "  +2.28%  "

$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.21"
$ws.Range("E19").Value = "  +6.27%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.51"
$ws.Range("E20").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "450.85"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.21"
$ws.Range("E22").Value = "  -2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  -1.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.45"
$ws.Range("E24").Value = "  +5.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.32"
$ws.Range("E25").Value = "  +9.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.27"
$ws.Range("E26").Value = "  +12.04%  "

$ws.Range("D27").Value = "3.452.56"
$ws.Range("E27").Value = "  +2.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "78.32"
$ws.Range("E28").Value = "  +0.56%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000128"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("B30").Value = "Dai"
$ws.Range("C30").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("E32").Value = "  -0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "573.57"
$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.155"
$ws.Range("E34").Value = "  +25.93%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.153"
$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.33"
$ws.Range("E38").Value = "  +0.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.28"
$ws.Range("E39").Value = "  +9.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.12"
$ws.Range("E42").Value = "  +14.26%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.05"
$ws.Range("E43").Value = "  +12.40%  "

$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.94"
$ws.Range("E44").Value = "  +3.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.70"
$ws.Range("E45").Value = "  -2.29%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "191.80"
$ws.Range("E47").Value = "  -0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "44.98"
$ws.Range("E48").Value = "  +4.03%  "

$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.783"
$ws.Range("E50").Value = "  -2.33%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.46"
$ws.Range("E51").Value = "  +1.85%  "
